# Momentum summary update: add "bias" columns (D, G) and a second
# "switched" breakdown column (L); recolor the table by section;
# update formulas/values to match the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Insert two new columns: one after C (new D), one after old F
#    (which, after the first insert, sits at column G).
# ---------------------------------------------------------------
$ws.Columns.Item(4).Insert()   # new column D: "Positive Momentum Bias "
$ws.Columns.Item(7).Insert()   # new column G: "Negative Momentum Bias"

# ---------------------------------------------------------------
# 2. Header row (row 1)
# ---------------------------------------------------------------
$ws.Range("D1").Value = "Positive Momentum Bias "
$ws.Range("G1").Value = "Negative Momentum Bias"
$ws.Range("K1").Value = "Momentum Switched & Positive Winners"
$ws.Range("L1").Value = "Momentum Switched & Negative Winners"

# ---------------------------------------------------------------
# 3. Data rows 2-6 for the two new "bias" columns
# ---------------------------------------------------------------
$ws.Range("D2").Value = 4.96
$ws.Range("D3").Value = 5.19
$ws.Range("D4").Value = 5.41
$ws.Range("D5").Value = 5.0999999999999996
$ws.Range("D6").Value = 5.04

$ws.Range("G2").Value = 4.9400000000000004
$ws.Range("G3").Value = 5.76
$ws.Range("G4").Value = 6.04
$ws.Range("G5").Value = 5.25
$ws.Range("G6").Value = 5.6

# New column L data (rows 2-6) - "Momentum Switched & Negative Winners"
$ws.Range("L2").Value = 3
$ws.Range("L3").Value = 2
$ws.Range("L4").Value = 5
$ws.Range("L5").Value = 6
$ws.Range("L6").Value = 8

# ---------------------------------------------------------------
# 4. Row 7 (totals) - add SUM formulas for new columns D, G, L and
#    extend the shared SUM formula range from F7:I7 to H7:K7
# ---------------------------------------------------------------
$ws.Range("D7").Formula = "=SUM(D2:D6)"
$ws.Range("G7").Formula = "=SUM(G2:G6)"
$ws.Range("L7").Formula = "=SUM(L2:L6)"

# ---------------------------------------------------------------
# 5. Row 8 (percent / summary row)
# ---------------------------------------------------------------
$ws.Range("B8").Formula = "=(B7+D7) / (B7 + C7)"
$ws.Range("C8").Formula = "=(C7-D7) / (C7 + B7)"
$ws.Range("D8").Value = 25.7
$ws.Range("E8").Formula = "=(E7+G7) / (E7 + F7)"
$ws.Range("F8").Formula = "=(F7-G7) / (F7 + E7)"
$ws.Range("G8").Value = 57
$ws.Range("K8").Formula = "=K7 / (J7)"
$ws.Range("L8").Formula = "=L7 / (J7)"

# ---------------------------------------------------------------
# 6. Fills: green for the "Wins/Bias" block (B:G), yellow stays on
#    the original summary block (now H:L), blue for the A "year"
#    column.
# ---------------------------------------------------------------
$yellow = 65535       # RGB(255,255,0)
$green  = 5287310     # RGB(80,208,146) -> actually set explicitly below
$blue   = 15773696    # RGB(0,176,240)

# Precise RGB() reimplementation (R + G*256 + B*65536), since RGB() isn't
# available as a function in this host:
# Yellow FFFF00 -> R255 G255 B0
$colYellow = 255 + (255 * 256) + (0 * 65536)
# Green 92D050 -> R146 G208 B80
$colGreen  = 146 + (208 * 256) + (80 * 65536)
# Blue 00B0F0 -> R0 G176 B240
$colBlue   = 0 + (176 * 256) + (240 * 65536)

$ws.Range("A1:A8").Interior.Color = $colBlue
$ws.Range("B1:G8").Interior.Color = $colGreen
$ws.Range("H1:L8").Interior.Color = $colYellow

# ---------------------------------------------------------------
# 7. Number formats: percent columns keep their "0%"-style format
#    (numFmtId 9) on the new bias-derived percent cells too.
# ---------------------------------------------------------------
$ws.Range("I8").NumberFormat = "0%"
$ws.Range("K8").NumberFormat = "0%"
$ws.Range("L8").NumberFormat = "0%"

# ---------------------------------------------------------------
# 8. Sheet view tweaks
# ---------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A7").Select()
